$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sema3a"
$ws.Range("C2").Value = "Nrp2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.5674196666666667
$ws.Range("H2").Value = 1.702259
$ws.Range("I2").Value = 0.07864125446886469
$ws.Range("J2").Value = 0.07864125446886468
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 34.52052333333333
$ws.Range("N2").Value = 103.56157
$ws.Range("O2").Value = 0.7684334662422598
$ws.Range("P2").Value = 0.7684334662422598
$ws.Range("Q2").Value = 19.58762384295889
$ws.Range("R2").Value = 176.28861458663
$ws.Range("S2").Value = 0.0604305717611493
$ws.Range("T2").Value = 0.06043057176114929

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sema3a"
$ws.Range("C3").Value = "Nrp2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.5674196666666667
$ws.Range("H3").Value = 1.702259
$ws.Range("I3").Value = 0.07864125446886469
$ws.Range("J3").Value = 0.07864125446886468
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.347618333333333
$ws.Range("N3").Value = 16.042855
$ws.Range("O3").Value = 0.1190390091234806
$ws.Range("P3").Value = 0.1190390091234805
$ws.Range("Q3").Value = 3.034343812160555
$ws.Range("R3").Value = 27.309094309445
$ws.Range("S3").Value = 0.00936137700820114
$ws.Range("T3").Value = 0.009361377008201138

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Sema3a"
$ws.Range("C4").Value = "Nrp2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.5674196666666667
$ws.Range("H4").Value = 1.702259
$ws.Range("I4").Value = 0.07864125446886469
$ws.Range("J4").Value = 0.07864125446886468
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.055101333333334
$ws.Range("N4").Value = 15.165304
$ws.Range("O4").Value = 0.1125275246342597
$ws.Range("P4").Value = 0.1125275246342597
$ws.Range("Q4").Value = 2.868363913526223
$ws.Range("R4").Value = 25.815275221736
$ws.Range("S4").Value = 0.008849305699514257
$ws.Range("T4").Value = 0.008849305699514256

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Sema3a"
$ws.Range("C5").Value = "Nrp2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7227763333333334
$ws.Range("H5").Value = 2.168329
$ws.Range("I5").Value = 0.1001728366019618
$ws.Range("J5").Value = 0.1001728366019618
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 34.52052333333333
$ws.Range("N5").Value = 103.56157
$ws.Range("O5").Value = 0.7684334662422598
$ws.Range("P5").Value = 0.7684334662422598
$ws.Range("Q5").Value = 24.95061727961444
$ws.Range("R5").Value = 224.55555551653
$ws.Range("S5").Value = 0.07697616005336502
$ws.Range("T5").Value = 0.07697616005336502

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Sema3a"
$ws.Range("C6").Value = "Nrp2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7227763333333334
$ws.Range("H6").Value = 2.168329
$ws.Range("I6").Value = 0.1001728366019618
$ws.Range("J6").Value = 0.1001728366019618
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 5.347618333333333
$ws.Range("N6").Value = 16.042855
$ws.Range("O6").Value = 0.1190390091234806
$ws.Range("P6").Value = 0.1190390091234805
$ws.Range("Q6").Value = 3.865131971032778
$ws.Range("R6").Value = 34.78618773929499
$ws.Range("S6").Value = 0.01192447521018586
$ws.Range("T6").Value = 0.01192447521018586

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sema3a"
$ws.Range("C7").Value = "Nrp2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.7227763333333334
$ws.Range("H7").Value = 2.168329
$ws.Range("I7").Value = 0.1001728366019618
$ws.Range("J7").Value = 0.1001728366019618
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.055101333333334
$ws.Range("N7").Value = 15.165304
$ws.Range("O7").Value = 0.1125275246342597
$ws.Range("P7").Value = 0.1125275246342597
$ws.Range("Q7").Value = 3.653707606335112
$ws.Range("R7").Value = 32.883368457016
$ws.Range("S7").Value = 0.01127220133841093
$ws.Range("T7").Value = 0.01127220133841093

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Sema3a"
$ws.Range("C8").Value = "Nrp2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.925096666666668
$ws.Range("H8").Value = 17.77529
$ws.Range("I8").Value = 0.8211859089291735
$ws.Range("J8").Value = 0.8211859089291734
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 34.52052333333333
$ws.Range("N8").Value = 103.56157
$ws.Range("O8").Value = 0.7684334662422598
$ws.Range("P8").Value = 0.7684334662422598
$ws.Range("Q8").Value = 204.5374377339222
$ws.Range("R8").Value = 1840.8369396053
$ws.Range("S8").Value = 0.6310267344277455
$ws.Range("T8").Value = 0.6310267344277455

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Sema3a"
$ws.Range("C9").Value = "Nrp2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.925096666666668
$ws.Range("H9").Value = 17.77529
$ws.Range("I9").Value = 0.8211859089291735
$ws.Range("J9").Value = 0.8211859089291734
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 5.347618333333333
$ws.Range("N9").Value = 16.042855
$ws.Range("O9").Value = 0.1190390091234806
$ws.Range("P9").Value = 0.1190390091234805
$ws.Range("Q9").Value = 31.68515556143889
$ws.Range("R9").Value = 285.16640005295
$ws.Range("S9").Value = 0.09775315690509356
$ws.Range("T9").Value = 0.09775315690509354

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Sema3a"
$ws.Range("C10").Value = "Nrp2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 5.925096666666668
$ws.Range("H10").Value = 17.77529
$ws.Range("I10").Value = 0.8211859089291735
$ws.Range("J10").Value = 0.8211859089291734
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.055101333333334
$ws.Range("N10").Value = 15.165304
$ws.Range("O10").Value = 0.1125275246342597
$ws.Range("P10").Value = 0.1125275246342597
$ws.Range("Q10").Value = 29.95196405979556
$ws.Range("R10").Value = 269.56767653816
$ws.Range("S10").Value = 0.09240601759633452
$ws.Range("T10").Value = 0.0924060175963345

